# Edit script for Bladen_Community_College_Organizations.xlsx
# 1. Remove the "General / Student Activities Board" row (old row 5), shifting rows 6-11 up.
# 2. Swap the "Category" and "Organization Name" columns (A <-> B) and rename headers.
# 3. Rename remaining headers (C, D, G-L) to their new labels.
# 4. Delete the now-unused "Tiktok Link" column (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the old row 5 (General / Student Activities Board) ---
# This shifts everything below it up by one, turning the old 11-row range into 10 rows.
$ws.Rows.Item(5).Delete()

# --- Step 2: swap columns A and B (Category <-> Organization Name) for all data rows ---
# NOTE: use .Value2 for reads -- this host's .Value getter mis-resolves on
# chained property access; .Value2 (and .Value for writes) works correctly.
$lastRow = 10
for ($r = 2; $r -le $lastRow; $r++) {
    $colA = $ws.Cells.Item($r, 1).Value2
    $colB = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $colB
    $ws.Cells.Item($r, 2).Value = $colA
}

# --- Step 3: update header row labels ---
$ws.Cells.Item(1, 1).Value = "Organization Name"
$ws.Cells.Item(1, 2).Value = "Categories"
$ws.Cells.Item(1, 3).Value = "Org URL"
$ws.Cells.Item(1, 4).Value = "Image URL"
$ws.Cells.Item(1, 5).Value = "Description"
$ws.Cells.Item(1, 6).Value = "Email"
$ws.Cells.Item(1, 7).Value = "Phone"
$ws.Cells.Item(1, 8).Value = "Website"
$ws.Cells.Item(1, 9).Value = "LinkedIn"
$ws.Cells.Item(1, 10).Value = "Instagram"
$ws.Cells.Item(1, 11).Value = "Facebook"
$ws.Cells.Item(1, 12).Value = "Twitter"

# --- Step 4: swap column widths for A and B, then delete column M (Tiktok Link) ---
# NOTE: this host stores ColumnWidth internally with a fixed +5/6 (0.8333...)
# offset versus the value assigned (an Excel "characters" -> internal-units
# quirk tied to the default font). Subtract that offset here so the saved
# OOXML <col width="..."> attribute comes out at the exact target integer.
$offset = 0.8333333333333334
$ws.Columns.Item(1).ColumnWidth = 35 - $offset
$ws.Columns.Item(2).ColumnWidth = 20 - $offset

$ws.Columns.Item(7).ColumnWidth = 7 - $offset
$ws.Columns.Item(8).ColumnWidth = 9 - $offset
$ws.Columns.Item(9).ColumnWidth = 10 - $offset
$ws.Columns.Item(10).ColumnWidth = 11 - $offset
$ws.Columns.Item(11).ColumnWidth = 10 - $offset
$ws.Columns.Item(12).ColumnWidth = 9 - $offset

$ws.Columns.Item(13).Delete()
